$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3582
$ws1.Range("G4").Value = 75
$ws1.Range("F5").Value = 3582
$ws1.Range("F6").Value = 259
$ws1.Range("F7").Value = 5111
$ws1.Range("F8").Value = 5111
$ws1.Range("F9").Value = 533
$ws1.Range("F10").Value = 360
$ws1.Range("F12").Value = 692
$ws1.Range("F14").Value = 93
$ws1.Range("F15").Value = 34
$ws1.Range("F17").Value = 317
$ws1.Range("F18").Value = 36
$ws1.Range("F20").Value = 156
$ws1.Range("F21").Value = 4
$ws1.Range("F23").Value = 4914
$ws1.Range("F24").Value = 4914
$ws1.Range("F26").Value = 41
$ws1.Range("F28").Value = 6044
$ws1.Range("F31").Value = 3222
$ws1.Range("F32").Value = 343
$ws1.Range("F33").Value = 713
$ws1.Range("F38").Value = 1024
$ws1.Range("F42").Value = 873
$ws1.Range("F43").Value = 998
$ws1.Range("F44").Value = 2029

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1121

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1121
$ws4.Range("F7").Value = 3582
$ws4.Range("G7").Value = 75
$ws4.Range("F8").Value = 3582
$ws4.Range("F9").Value = 259
$ws4.Range("F10").Value = 5111
$ws4.Range("F11").Value = 5111
$ws4.Range("F12").Value = 533
$ws4.Range("F13").Value = 360
$ws4.Range("F15").Value = 692
$ws4.Range("F17").Value = 93
$ws4.Range("F18").Value = 34
$ws4.Range("F20").Value = 317
$ws4.Range("F21").Value = 36
$ws4.Range("F24").Value = 156
$ws4.Range("F25").Value = 4
$ws4.Range("F27").Value = 4914
$ws4.Range("F28").Value = 4914
$ws4.Range("F30").Value = 41
$ws4.Range("F32").Value = 6044
$ws4.Range("F35").Value = 3222
$ws4.Range("F36").Value = 343
$ws4.Range("F37").Value = 713
$ws4.Range("F43").Value = 1024
$ws4.Range("F47").Value = 873
$ws4.Range("F48").Value = 998
$ws4.Range("F50").Value = 2029

Write-Host "Applied all updates"